$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update weekly triaged issues (Closed Issues = column B, Opened Issues = column C)
$ws.Range("B8").Value = 41
$ws.Range("C8").Value = 20

$ws.Range("B9").Value = 47
$ws.Range("C9").Value = 24

$ws.Range("B10").Value = 35
$ws.Range("C10").Value = 18

$ws.Range("B11").Value = 38
$ws.Range("C11").Value = 18

$ws.Range("B12").Value = 29
$ws.Range("C12").Value = 19

$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 24
